$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 23 new rows above the old row 777 ("Remessa" order block), ---
# --- shifting the old rows 777:786 down to 800:809.                    ---
$ws.Rows("777:799").Insert()

# Column A values for the new rows 777:799 (order/lot numbers). These are
# purely-numeric-looking strings that must be stored as TEXT (matching the
# rest of column A), not auto-converted to numbers. We stage each value in
# a scratch cell with a leading apostrophe (forcing text), then copy/paste
# *values only* into the destination - the destination cell already carries
# the sheet's normal style (s=3) from the row Insert above, and
# PasteSpecial(values) keeps that destination style instead of bringing the
# scratch cell's "quote prefix" style along with it.
$aValues = @(
    "80267767",
    "80267768",
    "80267768",
    "80267768",
    "80267768",
    "80267768",
    "80267768",
    "80267768",
    "80267768",
    "80267769",
    "80267769",
    "80267769",
    "80267771",
    "80267774",
    "80267776",
    "80267778",
    "80267779",
    "80267780",
    "80267781",
    "80267781",
    "80267782",
    "80267782",
    "80267784"
)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = 777 + $i
    $ws.Range("ZZ1").Value = "'" + $aValues[$i]
    $ws.Range("ZZ1").Copy()
    $ws.Range("A$row").PasteSpecial(-4163)  # xlPasteValues
}
$ws.Range("ZZ1").Clear()

# Column B (SKU codes) - already non-numeric text, assign directly.
$ws.Range("B777").Value = "20964-CTY-I"
$ws.Range("B778").Value = "10376-ARI-I"
$ws.Range("B779").Value = "10388-ARI-I"
$ws.Range("B780").Value = "10369-ARI-I"
$ws.Range("B781").Value = "10355-ARI-I"
$ws.Range("B782").Value = "10354-ARI-I"
$ws.Range("B783").Value = "10498-ARI-I"
$ws.Range("B784").Value = "10403-ARI-I"
$ws.Range("B785").Value = "10497-ARI-I"
$ws.Range("B786").Value = "26489-YAG-I"
$ws.Range("B787").Value = "11425-ROY-I"
$ws.Range("B788").Value = "10742-ROY-I"
$ws.Range("B789").Value = "10190-VIS-I"
$ws.Range("B790").Value = "10378-ARI-I"
$ws.Range("B791").Value = "10526-ARI-I"
$ws.Range("B792").Value = "10045-ARI-I"
$ws.Range("B793").Value = "13972-TDK-N"
$ws.Range("B794").Value = "10040-ARI-I"
$ws.Range("B795").Value = "10078-BLB-I"
$ws.Range("B796").Value = "10077-BLB-I"
$ws.Range("B797").Value = "10485-ARI-I"
$ws.Range("B798").Value = "10480-ARI-I"
$ws.Range("B799").Value = "23422-GPB-I"

# Column C (quantities).
$ws.Range("C777").Value = 1
$ws.Range("C778").Value = 2
$ws.Range("C779").Value = 0
$ws.Range("C780").Value = 0
$ws.Range("C781").Value = 2
$ws.Range("C782").Value = 2
$ws.Range("C783").Value = 2
$ws.Range("C784").Value = 2
$ws.Range("C785").Value = 1
$ws.Range("C786").Value = 12000
$ws.Range("C787").Value = 100000
$ws.Range("C788").Value = 45000
$ws.Range("C789").Value = 500
$ws.Range("C790").Value = 1
$ws.Range("C791").Value = 1
$ws.Range("C792").Value = 1
$ws.Range("C793").Value = 1000
$ws.Range("C794").Value = 1
$ws.Range("C795").Value = 8
$ws.Range("C796").Value = 35
$ws.Range("C797").Value = 1
$ws.Range("C798").Value = 1
$ws.Range("C799").Value = 100

# Match the final on-save cursor position recorded in the workbook.
$ws.Range("F8").Select() | Out-Null
